$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3935.5862
$ws.Range("I15").Value = 3935.5862
$ws.Range("K15").Value = 11806.7586
$ws.Range("M15").Value = -11637.7586
$ws.Range("H33").Value = 1330.8064
$ws.Range("I33").Value = 954.48
$ws.Range("K33").Value = 954.48
$ws.Range("M33").Value = -725.48
$ws.Range("H70").Value = 3479.1428
$ws.Range("J70").Value = 3889.5
$ws.Range("L70").Value = 11668.5
$ws.Range("N70").Value = -12208.5
$ws.Range("H73").Value = 3479.1428
$ws.Range("J73").Value = 3889.5
$ws.Range("L73").Value = 11668.5
$ws.Range("N73").Value = -13540.5
$ws.Range("H86").Value = 6649.75
$ws.Range("I86").Value = 7099.7334
$ws.Range("J86").Value = 5299.8
$ws.Range("K86").Value = 7099.7334
$ws.Range("L86").Value = 5299.8
$ws.Range("M86").Value = -5976.7334
$ws.Range("N86").Value = -7545.8
$ws.Range("H89").Value = 6649.75
$ws.Range("I89").Value = 7099.7334
$ws.Range("J89").Value = 5299.8
$ws.Range("K89").Value = 35498.667
$ws.Range("L89").Value = 26499
$ws.Range("M89").Value = -29882.667
$ws.Range("N89").Value = -37731
$ws.Range("H138").Value = 2614.0488
$ws.Range("I138").Value = 969.875
$ws.Range("J138").Value = 3012.6365
$ws.Range("K138").Value = 2909.625
$ws.Range("L138").Value = 9037.9095
$ws.Range("M138").Value = 2230.375
$ws.Range("N138").Value = -19317.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 829.4
$ws.Range("I4").Value = 287.25
$ws.Range("K4").Value = 287.25
$ws.Range("M4").Value = -171.25
$ws.Range("H32").Value = 8478805
$ws.Range("I32").Value = 8478805
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 8478805
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -8478518
$ws.Range("H61").Value = 29479048
$ws.Range("I61").Value = 83340740
$ws.Range("K61").Value = 83340740
$ws.Range("M61").Value = -83340528
$ws.Range("H97").Value = 638.5
$ws.Range("I97").Value = 514.7273
$ws.Range("K97").Value = 514.7273
$ws.Range("M97").Value = -18.72730000000001
$ws.Range("H122").Value = 1861
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2222
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 6666
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -11566
$ws.Range("H132").Value = 6187.08
$ws.Range("I132").Value = 3667.2778
$ws.Range("J132").Value = 12666.571
$ws.Range("K132").Value = 11001.8334
$ws.Range("L132").Value = 37999.713
$ws.Range("M132").Value = -8471.8334
$ws.Range("N132").Value = -43059.713
$ws.Range("H136").Value = 29479048
$ws.Range("I136").Value = 83340740
$ws.Range("K136").Value = 250022220
$ws.Range("M136").Value = -250019670

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -14887
$ws.Range("H107").Value = 2284
$ws.Range("I107").Value = 1926.25
$ws.Range("K107").Value = 1926.25
$ws.Range("M107").Value = -6.25
$ws.Range("H134").Value = 78518.38
$ws.Range("I134").Value = 1727.0834
$ws.Range("K134").Value = 5181.2502
$ws.Range("M134").Value = -2646.2502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 39818
$ws.Range("I51").Value = 22571.285
$ws.Range("K51").Value = 22571.285
$ws.Range("M51").Value = -21835.285
$ws.Range("H61").Value = 39818
$ws.Range("I61").Value = 22571.285
$ws.Range("K61").Value = 22571.285
$ws.Range("M61").Value = -22223.285
$ws.Range("H107").Value = 2053.25
$ws.Range("I107").Value = 1350
$ws.Range("K107").Value = 1350
$ws.Range("M107").Value = 570
$ws.Range("H132").Value = 3070.2856
$ws.Range("I132").Value = 3070.2856
$ws.Range("K132").Value = 9210.856800000001
$ws.Range("M132").Value = -6680.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 745.1539
$ws.Range("I5").Value = 745.1539
$ws.Range("K5").Value = 2235.4617
$ws.Range("M5").Value = -2123.4617
$ws.Range("H23").Value = 970.4286
$ws.Range("I23").Value = 932.3333
$ws.Range("K23").Value = 2796.9999
$ws.Range("M23").Value = -2561.9999
$ws.Range("H59").Value = 1881.8182
$ws.Range("I59").Value = 1500
$ws.Range("J59").Value = 2550
$ws.Range("K59").Value = 4500
$ws.Range("L59").Value = 7650
$ws.Range("M59").Value = -3960
$ws.Range("N59").Value = -8730
$ws.Range("H132").Value = 2544.0588
$ws.Range("I132").Value = 2553.5
$ws.Range("K132").Value = 22981.5
$ws.Range("M132").Value = -20451.5
$ws.Range("H135").Value = 745.1539
$ws.Range("I135").Value = 745.1539
$ws.Range("K135").Value = 6706.3851
$ws.Range("M135").Value = -4171.3851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I102").Value = 3815.9285
$ws.Range("J102").Value = 2333.3333
$ws.Range("K102").Value = 3815.9285
$ws.Range("L102").Value = 2333.3333
$ws.Range("M102").Value = -2193.9285
$ws.Range("N102").Value = -5577.3333
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H132").Value = 166670340
$ws.Range("I132").Value = 200003520
$ws.Range("K132").Value = 600010560
$ws.Range("M132").Value = -600008030
$ws.Range("H135").Value = 105999
$ws.Range("J135").Value = 105999
$ws.Range("L135").Value = 105999
$ws.Range("N135").Value = -116139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5011.3335
$ws.Range("I22").Value = 4700.5
$ws.Range("J22").Value = 5260
$ws.Range("K22").Value = 4700.5
$ws.Range("L22").Value = 5260
$ws.Range("M22").Value = -4405.5
$ws.Range("N22").Value = -5850
$ws.Range("H27").Value = 5011.3335
$ws.Range("I27").Value = 4700.5
$ws.Range("J27").Value = 5260
$ws.Range("K27").Value = 4700.5
$ws.Range("L27").Value = 5260
$ws.Range("M27").Value = -4593.5
$ws.Range("N27").Value = -5474
$ws.Range("H40").Value = 3728.6875
$ws.Range("I40").Value = 2566.4
$ws.Range("J40").Value = 5665.8335
$ws.Range("K40").Value = 2566.4
$ws.Range("L40").Value = 5665.8335
$ws.Range("M40").Value = -2430.4
$ws.Range("N40").Value = -5937.8335
$ws.Range("H68").Value = 4038.8
$ws.Range("I68").Value = 3775.889
$ws.Range("J68").Value = 4433.1665
$ws.Range("K68").Value = 3775.889
$ws.Range("L68").Value = 4433.1665
$ws.Range("M68").Value = -3026.889
$ws.Range("N68").Value = -5931.1665
$ws.Range("H71").Value = 4038.8
$ws.Range("I71").Value = 3775.889
$ws.Range("J71").Value = 4433.1665
$ws.Range("K71").Value = 18879.445
$ws.Range("L71").Value = 22165.8325
$ws.Range("M71").Value = -15135.445
$ws.Range("N71").Value = -29653.8325
$ws.Range("H87").Value = 111000
$ws.Range("J87").Value = 122000
$ws.Range("L87").Value = 122000
$ws.Range("N87").Value = -124246
$ws.Range("H90").Value = 111000
$ws.Range("J90").Value = 122000
$ws.Range("L90").Value = 366000
$ws.Range("N90").Value = -377232
$ws.Range("H100").Value = 2560.6667
$ws.Range("I100").Value = 2127.8572
$ws.Range("K100").Value = 2127.8572
$ws.Range("M100").Value = -1586.8572
$ws.Range("H120").Value = 66590
$ws.Range("J120").Value = 66590
$ws.Range("L120").Value = 66590
$ws.Range("N120").Value = -76266
$ws.Range("H136").Value = 48650.16
$ws.Range("I136").Value = 2912.5
$ws.Range("J136").Value = 231600.8
$ws.Range("K136").Value = 8737.5
$ws.Range("L136").Value = 694802.3999999999
$ws.Range("M136").Value = -6187.5
$ws.Range("N136").Value = -699902.3999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1030
$ws.Range("I23").Value = 540
$ws.Range("K23").Value = 540
$ws.Range("M23").Value = -311
$ws.Range("H96").Value = 3797
$ws.Range("I96").Value = 1100
$ws.Range("K96").Value = 1100
$ws.Range("M96").Value = 273
$ws.Range("H107").Value = 17242180
$ws.Range("I107").Value = 27778784
$ws.Range("J107").Value = 465.45456
$ws.Range("K107").Value = 83336352
$ws.Range("L107").Value = 1396.36368
$ws.Range("M107").Value = -83334432
$ws.Range("N107").Value = -5236.36368
$ws.Range("H122").Value = 5288.28
$ws.Range("I122").Value = 4010.3333
$ws.Range("J122").Value = 6467.923
$ws.Range("K122").Value = 12030.9999
$ws.Range("L122").Value = 19403.769
$ws.Range("M122").Value = -9580.999899999999
$ws.Range("N122").Value = -24303.769
$ws.Range("H132").Value = 2250.3
$ws.Range("I132").Value = 1542.96
$ws.Range("K132").Value = 4628.88
$ws.Range("M132").Value = -2098.88
$ws.Range("H136").Value = 3299.6924
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100
